$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Timp2"
$ws.Cells.Item(2, 3).Value = "Itga3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 13.032878
$ws.Cells.Item(2, 8).Value = 39.098634
$ws.Cells.Item(2, 9).Value = 0.02949184097968156
$ws.Cells.Item(2, 10).Value = 0.02949184097968156
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 8.269168666666666
$ws.Cells.Item(2, 14).Value = 24.807506
$ws.Cells.Item(2, 15).Value = 0.671680253471746
$ws.Cells.Item(2, 16).Value = 0.671680253471746
$ws.Cells.Item(2, 17).Value = 107.7710663940893
$ws.Cells.Item(2, 18).Value = 969.939597546804
$ws.Cells.Item(2, 19).Value = 0.01980908722458094
$ws.Cells.Item(2, 20).Value = 0.01980908722458094

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Timp2"
$ws.Cells.Item(3, 3).Value = "Itga3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 13.032878
$ws.Cells.Item(3, 8).Value = 39.098634
$ws.Cells.Item(3, 9).Value = 0.02949184097968156
$ws.Cells.Item(3, 10).Value = 0.02949184097968156
$ws.Cells.Item(3, 11).Value = 2.0
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.121294
$ws.Cells.Item(3, 14).Value = 0.363882
$ws.Cells.Item(3, 15).Value = 0.009852354928133683
$ws.Cells.Item(3, 16).Value = 0.009852354928133683
$ws.Cells.Item(3, 17).Value = 1.580809904132
$ws.Cells.Item(3, 18).Value = 14.227289137188
$ws.Cells.Item(3, 19).Value = 0.0002905640848159006
$ws.Cells.Item(3, 20).Value = 0.0002905640848159006

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Timp2"
$ws.Cells.Item(4, 3).Value = "Itga3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 13.032878
$ws.Cells.Item(4, 8).Value = 39.098634
$ws.Cells.Item(4, 9).Value = 0.02949184097968156
$ws.Cells.Item(4, 10).Value = 0.02949184097968156
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 3.920705666666667
$ws.Cells.Item(4, 14).Value = 11.762117
$ws.Cells.Item(4, 15).Value = 0.3184673916001203
$ws.Cells.Item(4, 16).Value = 0.3184673916001203
$ws.Cells.Item(4, 17).Value = 51.09807862757534
$ws.Cells.Item(4, 18).Value = 459.882707648178
$ws.Cells.Item(4, 19).Value = 0.009392189670284724
$ws.Cells.Item(4, 20).Value = 0.009392189670284724

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Timp2"
$ws.Cells.Item(5, 3).Value = "Itga3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 395.9197996666667
$ws.Cells.Item(5, 8).Value = 1187.759399
$ws.Cells.Item(5, 9).Value = 0.8959190573622122
$ws.Cells.Item(5, 10).Value = 0.8959190573622122
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 8.269168666666666
$ws.Cells.Item(5, 14).Value = 24.807506
$ws.Cells.Item(5, 15).Value = 0.671680253471746
$ws.Cells.Item(5, 16).Value = 0.671680253471746
$ws.Cells.Item(5, 17).Value = 3273.927601916544
$ws.Cells.Item(5, 18).Value = 29465.34841724889
$ws.Cells.Item(5, 19).Value = 0.6017711395392185
$ws.Cells.Item(5, 20).Value = 0.6017711395392185

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Timp2"
$ws.Cells.Item(6, 3).Value = "Itga3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 395.9197996666667
$ws.Cells.Item(6, 8).Value = 1187.759399
$ws.Cells.Item(6, 9).Value = 0.8959190573622122
$ws.Cells.Item(6, 10).Value = 0.8959190573622122
$ws.Cells.Item(6, 11).Value = 2.0
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.121294
$ws.Cells.Item(6, 14).Value = 0.363882
$ws.Cells.Item(6, 15).Value = 0.009852354928133683
$ws.Cells.Item(6, 16).Value = 0.009852354928133683
$ws.Cells.Item(6, 17).Value = 48.02269618076867
$ws.Cells.Item(6, 18).Value = 432.204265626918
$ws.Cells.Item(6, 19).Value = 0.008826912540011476
$ws.Cells.Item(6, 20).Value = 0.008826912540011476

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Timp2"
$ws.Cells.Item(7, 3).Value = "Itga3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 395.9197996666667
$ws.Cells.Item(7, 8).Value = 1187.759399
$ws.Cells.Item(7, 9).Value = 0.8959190573622122
$ws.Cells.Item(7, 10).Value = 0.8959190573622122
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 3.920705666666667
$ws.Cells.Item(7, 14).Value = 11.762117
$ws.Cells.Item(7, 15).Value = 0.3184673916001203
$ws.Cells.Item(7, 16).Value = 0.3184673916001203
$ws.Cells.Item(7, 17).Value = 1552.285002098631
$ws.Cells.Item(7, 18).Value = 13970.56501888768
$ws.Cells.Item(7, 19).Value = 0.2853210052829823
$ws.Cells.Item(7, 20).Value = 0.2853210052829823

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Timp2"
$ws.Cells.Item(8, 3).Value = "Itga3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 32.96202033333334
$ws.Cells.Item(8, 8).Value = 98.886061
$ws.Cells.Item(8, 9).Value = 0.07458910165810628
$ws.Cells.Item(8, 10).Value = 0.07458910165810628
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 8.269168666666666
$ws.Cells.Item(8, 14).Value = 24.807506
$ws.Cells.Item(8, 15).Value = 0.671680253471746
$ws.Cells.Item(8, 16).Value = 0.671680253471746
$ws.Cells.Item(8, 17).Value = 272.5685057304295
$ws.Cells.Item(8, 18).Value = 2453.116551573866
$ws.Cells.Item(8, 19).Value = 0.05010002670794666
$ws.Cells.Item(8, 20).Value = 0.05010002670794666

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Timp2"
$ws.Cells.Item(9, 3).Value = "Itga3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 32.96202033333334
$ws.Cells.Item(9, 8).Value = 98.886061
$ws.Cells.Item(9, 9).Value = 0.07458910165810628
$ws.Cells.Item(9, 10).Value = 0.07458910165810628
$ws.Cells.Item(9, 11).Value = 2.0
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.121294
$ws.Cells.Item(9, 14).Value = 0.363882
$ws.Cells.Item(9, 15).Value = 0.009852354928133683
$ws.Cells.Item(9, 16).Value = 0.009852354928133683
$ws.Cells.Item(9, 17).Value = 3.998095294311334
$ws.Cells.Item(9, 18).Value = 35.982857648802
$ws.Cells.Item(9, 19).Value = 0.0007348783033063077
$ws.Cells.Item(9, 20).Value = 0.0007348783033063077

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Timp2"
$ws.Cells.Item(10, 3).Value = "Itga3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 32.96202033333334
$ws.Cells.Item(10, 8).Value = 98.886061
$ws.Cells.Item(10, 9).Value = 0.07458910165810628
$ws.Cells.Item(10, 10).Value = 0.07458910165810628
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 3.920705666666667
$ws.Cells.Item(10, 14).Value = 11.762117
$ws.Cells.Item(10, 15).Value = 0.3184673916001203
$ws.Cells.Item(10, 16).Value = 0.3184673916001203
$ws.Cells.Item(10, 17).Value = 129.2343799056819
$ws.Cells.Item(10, 18).Value = 1163.109419151137
$ws.Cells.Item(10, 19).Value = 0.02375419664685332
$ws.Cells.Item(10, 20).Value = 0.02375419664685332
